$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "model best score"

$ws.Range("B2").Value = 0.88175000000000003
$ws.Range("B3").Value = 0.81933999999999996
$ws.Range("B5").Value = 0.79074

$ws.Range("A8").Value = "SVM"
$ws.Range("B8").Value = 0.79200999999999999

$ws.Range("A9").Value = "NN"
$ws.Range("B9").Value = 0.79871000000000003

$ws.Range("D11").Select()
